$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new price records were added for Femacal de La Calera - Zapallo.
# They belong right before the existing row 294, so push the existing
# rows 294:338 down by two (to 296:340) and populate the freed rows.
$ws.Rows("294:295").Insert()

# New row 294: Camote, 1a (guarda)
$ws.Cells.Item(294, 1).Value = 3
$ws.Cells.Item(294, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(294, 3).Value = "Coquimbo"
$ws.Cells.Item(294, 4).Value = 44491
$ws.Cells.Item(294, 5).Value = 5
$ws.Cells.Item(294, 6).Value = 100112045
$ws.Cells.Item(294, 7).Value = "Zapallo"
$ws.Cells.Item(294, 8).Value = "Camote"
$ws.Cells.Item(294, 9).Value = "1a (guarda)"
$ws.Cells.Item(294, 10).Value = 310
$ws.Cells.Item(294, 11).Value = 600
$ws.Cells.Item(294, 12).Value = 700
$ws.Cells.Item(294, 13).Value = 648
$ws.Cells.Item(294, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(294, 15).Value = "Provincia de Talca"
$ws.Cells.Item(294, 16).Value = 648
$ws.Cells.Item(294, 17).Value = 1
$ws.Cells.Item(294, 18).Value = "Hortaliza"

# New row 295: Paine, 1a (guarda)
$ws.Cells.Item(295, 1).Value = 3
$ws.Cells.Item(295, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(295, 3).Value = "Coquimbo"
$ws.Cells.Item(295, 4).Value = 44491
$ws.Cells.Item(295, 5).Value = 5
$ws.Cells.Item(295, 6).Value = 100112045
$ws.Cells.Item(295, 7).Value = "Zapallo"
$ws.Cells.Item(295, 8).Value = "Paine"
$ws.Cells.Item(295, 9).Value = "1a (guarda)"
$ws.Cells.Item(295, 10).Value = 150
$ws.Cells.Item(295, 11).Value = 400
$ws.Cells.Item(295, 12).Value = 400
$ws.Cells.Item(295, 13).Value = 400
$ws.Cells.Item(295, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(295, 15).Value = "Provincia de Talca"
$ws.Cells.Item(295, 16).Value = 400
$ws.Cells.Item(295, 17).Value = 1
$ws.Cells.Item(295, 18).Value = "Hortaliza"
